$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B. This shifts the existing headers
# (Enrollment No., Department, Name, DOB, Admission Year, Phone Number,
# Auto Generated Password) from B..H to C..I.
$ws.Range("B1").EntireColumn.Insert()

# C1 currently holds the old "Enrollment No." text (shifted from old B1).
# Overwrite it with the new header text.
$ws.Range("C1").Value = "University Roll Number."

# Set the new B1 header text.
$ws.Range("B1").Value = "University Registration Number"

# Copy C1's formatting (bold header style) onto the new B1 cell.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)

# Resize the columns to match the new layout.
$ws.Columns.Item(2).ColumnWidth = 27
$ws.Columns.Item(3).ColumnWidth = 26.333333333333336
$ws.Columns.Item(4).ColumnWidth = 18
$ws.Columns.Item(5).ColumnWidth = 16.333333333333336
$ws.Columns.Item(6).ColumnWidth = 15.166666666666668
$ws.Columns.Item(7).ColumnWidth = 18.5
$ws.Columns.Item(8).ColumnWidth = 20.5
$ws.Columns.Item(9).ColumnWidth = 30.5

# Update the selection to B1 (matches the saved view state in the workbook).
$null = $ws.Range("B1").Select()
